# This script updates the "Price" (column D) and "Volume(1h)" (column E)
# values on the "cryptos" worksheet to reflect refreshed market data, as
# produced by the scheduled GitHub Actions symbol-list update.
# All target cells are plain text cells (t="inlineStr" in the OOXML), so
# each value is written with an explicit Text number format to prevent
# Excel's automatic number/percentage conversion, and the cell style is
# then reset back to "Normal" so no stray formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates: cell address -> new text value
$updates = @{
    "D2" = "287.51"
    "E2" = "-0.81%"
    "D3" = "30.97"
    "E3" = "1.53%"
    "D4" = "4.930"
    "E4" = "-0.22%"
    "D5" = "0.07313"
    "E5" = "1.83%"
    "D6" = "2.370"
    "E6" = "30.13%"
    "D7" = "7.732"
    "E7" = "1.48%"
    "D8" = "0.9036"
    "E8" = "0.77%"
    "D9" = "0.09371"
    "E9" = "21.26%"
    "E10" = "2.17%"
    "D11" = "0.08184"
    "E11" = "2.99%"
    "D12" = "0.03127"
    "E12" = "2.89%"
    "D13" = "0.09930"
    "D14" = "0.001503"
    "E14" = "-0.10%"
    "D15" = "0.005774"
    "E15" = "0.68%"
    "D16" = "3.496"
    "E16" = "0.90%"
    "D17" = "3.723"
    "E17" = "-0.36%"
    "D18" = "2.080"
    "E18" = "0.14%"
    "D19" = "0.3328"
    "E19" = "0.28%"
    "D20" = "0.1331"
    "E20" = "4.09%"
    "D21" = "4.214"
    "E21" = "4.47%"
    "D22" = "0.2099"
    "E22" = "-12.13%"
    "D23" = "0.04512"
    "E23" = "0.13%"
    "D24" = "0.001210"
    "E24" = "-0.43%"
    "D25" = "0.004162"
    "E25" = "-9.88%"
    "D26" = "0.0001301"
    "E26" = "-0.07%"
    "D39" = "0.01575"
    "E39" = "0.75%"
    "D40" = "0.04445"
    "E40" = "2.21%"
    "D41" = "0.007376"
    "E41" = "1.03%"
    "D42" = "0.009491"
    "E42" = "-5.81%"
    "D43" = "0.1325"
    "E43" = "1.86%"
    "D44" = "0.002241"
    "E44" = "9.19%"
    "D45" = "0.008958"
    "E45" = "-5.74%"
    "D46" = "0.00006115"
    "E46" = "2.18%"
    "E47" = "-0.12%"
    "D48" = "2.507"
    "E48" = "8.79%"
    "E50" = "-0.12%"
    "E51" = "-0.12%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.Style = "Normal"
}

Write-Output ("Updated " + $updates.Count + " cells")
